$wb = $excel.ActiveWorkbook

$tc01 = $wb.Worksheets.Item("TC01")

# Duplicate the TC01 sheet (keeps formatting/merged cells identical) and
# place the copy right after TC01, before "base form".
$tc01.Copy([System.Reflection.Missing]::Value, $tc01)

$tc02 = $wb.Worksheets.Item("TC01 (2)")
$tc02.Name = "TC02"

# Update the new sheet's header cells for the logout test case.
$tc02.Range("B1").Value = "Logout endpoint"
$tc02.Range("A1").Value = "TC02"

# Make the new sheet the active / selected tab.
$tc02.Activate()
$tc02.Select()
